# Refresh the "cryptos" price/volume table (GitHub Actions scheduled update).
# Column D ("Price") values look numeric (e.g. "72.617.97", "0.730") but are
# stored as TEXT in the workbook (thousands separated by "." rather than ",",
# and trailing zeros that must be preserved). Setting NumberFormat to "@"
# (Text) before assigning the Value keeps Excel from re-interpreting the
# string as a number/date, and resetting the Style to "Normal" afterwards
# avoids leaving a stray text-format style on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.617.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.049.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.63%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.730"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +20.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.040.70"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.57%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.177"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.695.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.060.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.02%  "
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.134"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.486.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.89"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "105.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +20.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("E31").Value = "  +16.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.72%  "
$ws.Range("E33").Value = "  +4.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "679.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "67.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.60%  "
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.434"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0867"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "  +4.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.158"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.08%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +16.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000274"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.18%  "
